$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data for June 2024 baseline

# Row 11: Tax Wedge for Business Assets
$ws.Range("B11").Value = 1.1
$ws.Range("C11").Value = 1.15
$ws.Range("D11").Value = 1.26
$ws.Range("E11").Value = 1.32
$ws.Range("F11").Value = 1.32
$ws.Range("G11").Value = 1.32
$ws.Range("H11").Value = 1.33
$ws.Range("I11").Value = 1.33
$ws.Range("J11").Value = 1.34
$ws.Range("K11").Value = 1.34
$ws.Range("L11").Value = 1.34

# Row 14: Equity-financed
$ws.Range("B14").Value = 1.38
$ws.Range("C14").Value = 1.43
$ws.Range("D14").Value = 1.6
$ws.Range("E14").Value = 1.66
$ws.Range("F14").Value = 1.66
$ws.Range("G14").Value = 1.66
$ws.Range("H14").Value = 1.67
$ws.Range("I14").Value = 1.67
$ws.Range("J14").Value = 1.68
$ws.Range("K14").Value = 1.68
$ws.Range("L14").Value = 1.68

# Row 15: Debt-financed
$ws.Range("B15").Value = 0.52
$ws.Range("C15").Value = 0.55
$ws.Range("D15").Value = 0.54
$ws.Range("E15").Value = 0.6
$ws.Range("F15").Value = 0.6
$ws.Range("G15").Value = 0.6
$ws.Range("H15").Value = 0.61
$ws.Range("I15").Value = 0.61
$ws.Range("J15").Value = 0.61
$ws.Range("K15").Value = 0.61
$ws.Range("L15").Value = 0.62

# Row 17: Difference between sources of financing
$ws.Range("B17").Value = 0.86
$ws.Range("C17").Value = 0.88
$ws.Range("D17").Value = 1.06
$ws.Range("E17").Value = 1.07
$ws.Range("F17").Value = 1.07
$ws.Range("G17").Value = 1.07
$ws.Range("H17").Value = 1.06
$ws.Range("I17").Value = 1.06
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 1.07
$ws.Range("L17").Value = 1.07

# Row 20: C corporations
$ws.Range("B20").Value = 1.23
$ws.Range("C20").Value = 1.27
$ws.Range("D20").Value = 1.32
$ws.Range("E20").Value = 1.38
$ws.Range("F20").Value = 1.38
$ws.Range("G20").Value = 1.38
$ws.Range("H20").Value = 1.39
$ws.Range("I20").Value = 1.39
$ws.Range("J20").Value = 1.4
$ws.Range("K20").Value = 1.4
$ws.Range("L20").Value = 1.4

# Row 21: Pass-through entities
$ws.Range("B21").Value = 0.85
$ws.Range("C21").Value = 0.9
$ws.Range("D21").Value = 1.13
$ws.Range("E21").Value = 1.21
$ws.Range("F21").Value = 1.21
$ws.Range("G21").Value = 1.21
$ws.Range("H21").Value = 1.21
$ws.Range("I21").Value = 1.21
$ws.Range("J21").Value = 1.22
$ws.Range("K21").Value = 1.21
$ws.Range("L21").Value = 1.22

# Row 23: Difference between legal forms of organization
$ws.Range("B23").Value = 0.38
$ws.Range("C23").Value = 0.37
$ws.Range("D23").Value = 0.19
$ws.Range("E23").Value = 0.17
$ws.Range("F23").Value = 0.17
$ws.Range("G23").Value = 0.17
$ws.Range("H23").Value = 0.18
$ws.Range("I23").Value = 0.18
$ws.Range("J23").Value = 0.18
$ws.Range("K23").Value = 0.18
$ws.Range("L23").Value = 0.18

# Row 25: Weighted mean absolute difference between all asset pairs
$ws.Range("B25").Value = 0.85
$ws.Range("C25").Value = 0.8
$ws.Range("D25").Value = 0.79
$ws.Range("E25").Value = 0.74
$ws.Range("F25").Value = 0.74
$ws.Range("G25").Value = 0.73
$ws.Range("H25").Value = 0.73
$ws.Range("I25").Value = 0.73
$ws.Range("J25").Value = 0.73
$ws.Range("K25").Value = 0.73
$ws.Range("L25").Value = 0.73

# Row 26: Weighted mean absolute difference between all industry pairs
$ws.Range("B26").Value = 0.11
$ws.Range("C26").Value = 0.1
$ws.Range("D26").Value = 0.09
$ws.Range("E26").Value = 0.06
$ws.Range("F26").Value = 0.05
$ws.Range("G26").Value = 0.04
$ws.Range("H26").Value = 0.03
$ws.Range("I26").Value = 0.03
$ws.Range("J26").Value = 0.03
$ws.Range("K26").Value = 0.03
$ws.Range("L26").Value = 0.03

# Row 28: Tax Wedge for Owner-Occupied Housing
$ws.Range("B28").Value = 0.12
$ws.Range("C28").Value = 0.09
$ws.Range("D28").Value = -0.42
$ws.Range("E28").Value = -0.43
$ws.Range("F28").Value = -0.44
$ws.Range("G28").Value = -0.44
$ws.Range("H28").Value = -0.44
$ws.Range("I28").Value = -0.44
$ws.Range("J28").Value = -0.45
$ws.Range("K28").Value = -0.45
$ws.Range("L28").Value = -0.45

# Row 29: Difference between owner-occupied housing and business assets
$ws.Range("B29").Value = -0.99
$ws.Range("C29").Value = -1.06
$ws.Range("D29").Value = -1.68
$ws.Range("E29").Value = -1.75
$ws.Range("F29").Value = -1.76
$ws.Range("G29").Value = -1.76
$ws.Range("H29").Value = -1.77
$ws.Range("I29").Value = -1.77
$ws.Range("J29").Value = -1.78
$ws.Range("K29").Value = -1.78
$ws.Range("L29").Value = -1.79

# Row 32: Equity-financed
$ws.Range("B32").Value = -0.01
$ws.Range("C32").Value = -0.01
$ws.Range("D32").Value = -0.19
$ws.Range("E32").Value = -0.19
$ws.Range("F32").Value = -0.19
$ws.Range("G32").Value = -0.19
$ws.Range("H32").Value = -0.19
$ws.Range("I32").Value = -0.19
$ws.Range("J32").Value = -0.19
$ws.Range("K32").Value = -0.19
$ws.Range("L32").Value = -0.19

# Row 33: Debt-financed
$ws.Range("B33").Value = 0.34
$ws.Range("C33").Value = 0.27
$ws.Range("D33").Value = -0.83
$ws.Range("E33").Value = -0.86
$ws.Range("F33").Value = -0.87
$ws.Range("G33").Value = -0.88
$ws.Range("H33").Value = -0.88
$ws.Range("I33").Value = -0.88
$ws.Range("J33").Value = -0.89
$ws.Range("K33").Value = -0.89
$ws.Range("L33").Value = -0.89

# Row 35: Difference between sources of financing
$ws.Range("B35").Value = -0.35
$ws.Range("C35").Value = -0.28
$ws.Range("D35").Value = 0.64
$ws.Range("E35").Value = 0.66
$ws.Range("F35").Value = 0.68
$ws.Range("G35").Value = 0.69
$ws.Range("H35").Value = 0.68
$ws.Range("I35").Value = 0.69
$ws.Range("J35").Value = 0.7
$ws.Range("K35").Value = 0.7
$ws.Range("L35").Value = 0.7

# Row 38: Tax Wedge for Owner-Occupied Housing Structures
$ws.Range("B38").Value = 0.12
$ws.Range("C38").Value = 0.09
$ws.Range("D38").Value = -0.42
$ws.Range("E38").Value = -0.43
$ws.Range("F38").Value = -0.44
$ws.Range("G38").Value = -0.44
$ws.Range("H38").Value = -0.44
$ws.Range("I38").Value = -0.44
$ws.Range("J38").Value = -0.45
$ws.Range("K38").Value = -0.45
$ws.Range("L38").Value = -0.45

# Row 39: Tax Wedge for Renter-Occupied Housing Structures
$ws.Range("B39").Value = 1.45
$ws.Range("C39").Value = 1.44
$ws.Range("D39").Value = 1.51
$ws.Range("E39").Value = 1.5
$ws.Range("F39").Value = 1.5
$ws.Range("G39").Value = 1.5
$ws.Range("H39").Value = 1.5
$ws.Range("I39").Value = 1.5
$ws.Range("J39").Value = 1.5
$ws.Range("K39").Value = 1.51
$ws.Range("L39").Value = 1.51

# Row 41: Difference between owner- and renter-occupied housing structures
$ws.Range("B41").Value = -1.33
$ws.Range("C41").Value = -1.35
$ws.Range("D41").Value = -1.93
$ws.Range("E41").Value = -1.94
$ws.Range("F41").Value = -1.94
$ws.Range("G41").Value = -1.94
$ws.Range("H41").Value = -1.94
$ws.Range("I41").Value = -1.94
$ws.Range("J41").Value = -1.95
$ws.Range("K41").Value = -1.95
$ws.Range("L41").Value = -1.96
